$wb = $excel.ActiveWorkbook

# The "Users" sheet holds a single name in column A.
# Replace the old name with the new one.
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Drew Koecher"

# Select a cell on the Users sheet (as captured in the saved view state)
# and make it the active sheet/tab.
$usersSheet.Activate()
$usersSheet.Range("N18").Select()
